$wb = $excel.ActiveWorkbook

# Update market-board derived profit figures per scheduled data refresh.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1598.5
$ws.Range("I125").Value = 2032
$ws.Range("J125").Value = 1454
$ws.Range("K125").Value = 18288
$ws.Range("L125").Value = 13086
$ws.Range("M125").Value = -15828
$ws.Range("N125").Value = -18006

$ws.Range("H129").Value = 853.0244
$ws.Range("I129").Value = 820.9524
$ws.Range("J129").Value = 886.7
$ws.Range("K129").Value = 2462.8572
$ws.Range("L129").Value = 2660.1
$ws.Range("M129").Value = 2537.1428
$ws.Range("N129").Value = -12660.1

$ws.Range("H137").Value = 8319.054
$ws.Range("I137").Value = 21543.334
$ws.Range("J137").Value = 5759.516
$ws.Range("K137").Value = 64630.00199999999
$ws.Range("L137").Value = 17278.548
$ws.Range("M137").Value = -62080.00199999999
$ws.Range("N137").Value = -22378.548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2199.875
$ws.Range("I2").Value = 2199.875
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2199.875
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2086.875
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 19540.305
$ws.Range("I32").Value = 20897.371
$ws.Range("J32").Value = 15469.111
$ws.Range("K32").Value = 20897.371
$ws.Range("L32").Value = 15469.111
$ws.Range("M32").Value = -20610.371
$ws.Range("N32").Value = -16043.111

$ws.Range("H116").Value = 2199.875
$ws.Range("I116").Value = 2199.875
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2199.875
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 94.125
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 2121.2856
$ws.Range("I122").Value = 1789.4546
$ws.Range("J122").Value = 3338
$ws.Range("K122").Value = 5368.3638
$ws.Range("L122").Value = 10014
$ws.Range("M122").Value = -2918.3638
$ws.Range("N122").Value = -14914

$ws.Range("H132").Value = 13890892
$ws.Range("I132").Value = 16130730
$ws.Range("J132").Value = 3898.8
$ws.Range("K132").Value = 48392190
$ws.Range("L132").Value = 11696.4
$ws.Range("M132").Value = -48389660
$ws.Range("N132").Value = -16756.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2199.875
$ws.Range("I3").Value = 2199.875
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2199.875
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2085.875
$ws.Range("N3").ClearContents()

$ws.Range("H105").Value = 2384.3333
$ws.Range("I105").Value = 1806
$ws.Range("J105").Value = 7011
$ws.Range("K105").Value = 1806
$ws.Range("L105").Value = 7011
$ws.Range("M105").Value = -59
$ws.Range("N105").Value = -10505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3748.5
$ws.Range("I22").Value = 345.5
$ws.Range("J22").Value = 5450
$ws.Range("K22").Value = 345.5
$ws.Range("L22").Value = 5450
$ws.Range("M22").Value = 4.5
$ws.Range("N22").Value = -6150

$ws.Range("H99").Value = 1936.9565
$ws.Range("I99").Value = 2202
$ws.Range("J99").Value = 1821
$ws.Range("K99").Value = 2202
$ws.Range("L99").Value = 1821
$ws.Range("M99").Value = -704
$ws.Range("N99").Value = -4817

$ws.Range("H126").Value = 1936.9565
$ws.Range("I126").Value = 2202
$ws.Range("J126").Value = 1821
$ws.Range("K126").Value = 6606
$ws.Range("L126").Value = 5463
$ws.Range("M126").Value = -4136
$ws.Range("N126").Value = -10403

$ws.Range("H132").Value = 670859.9
$ws.Range("I132").Value = 2522.5454
$ws.Range("J132").Value = 1406030.9
$ws.Range("K132").Value = 7567.6362
$ws.Range("L132").Value = 4218092.699999999
$ws.Range("M132").Value = -5037.6362
$ws.Range("N132").Value = -4223152.699999999

$ws.Range("H134").Value = 637936.1
$ws.Range("I134").Value = 1373.2667
$ws.Range("J134").Value = 2001999.4
$ws.Range("K134").Value = 4119.800099999999
$ws.Range("L134").Value = 6005998.199999999
$ws.Range("M134").Value = -1584.800099999999
$ws.Range("N134").Value = -6011068.199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1610.3125
$ws.Range("I34").Value = 225.75
$ws.Range("J34").Value = 2071.8333
$ws.Range("K34").Value = 677.25
$ws.Range("L34").Value = 6215.499899999999
$ws.Range("M34").Value = -593.25
$ws.Range("N34").Value = -6383.499899999999

$ws.Range("H39").Value = 2509.0908
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2509.0908
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7527.2724
$ws.Range("N39").Value = -8115.2724

$ws.Range("H55").Value = 3000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354

$ws.Range("H59").Value = 1744.2858
$ws.Range("I59").Value = 1201.6666
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 3604.9998
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = -3064.9998
$ws.Range("N59").Value = -16080

$ws.Range("H64").Value = 4294.706
$ws.Range("I64").Value = 2400
$ws.Range("J64").Value = 5084.1665
$ws.Range("K64").Value = 7200
$ws.Range("L64").Value = 15252.4995
$ws.Range("M64").Value = -6930
$ws.Range("N64").Value = -15792.4995

$ws.Range("H67").Value = 4294.706
$ws.Range("I67").Value = 2400
$ws.Range("J67").Value = 5084.1665
$ws.Range("K67").Value = 7200
$ws.Range("L67").Value = 15252.4995
$ws.Range("M67").Value = -6264
$ws.Range("N67").Value = -17124.4995

$ws.Range("H68").Value = 11286108
$ws.Range("I68").Value = 6945436.5
$ws.Range("J68").Value = 15626780
$ws.Range("K68").Value = 20836309.5
$ws.Range("L68").Value = 46880340
$ws.Range("M68").Value = -20835498.5
$ws.Range("N68").Value = -46881962

$ws.Range("H71").Value = 11286108
$ws.Range("I71").Value = 6945436.5
$ws.Range("J71").Value = 15626780
$ws.Range("K71").Value = 62508928.5
$ws.Range("L71").Value = 140641020
$ws.Range("M71").Value = -62504872.5
$ws.Range("N71").Value = -140649132

$ws.Range("H97").Value = 51512.5
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 58828.57
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 176485.71
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -177477.71

$ws.Range("H121").Value = 258848.08
$ws.Range("I121").Value = 370
$ws.Range("J121").Value = 443475.28
$ws.Range("K121").Value = 1110
$ws.Range("L121").Value = 1330425.84
$ws.Range("M121").Value = 200
$ws.Range("N121").Value = -1333045.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H92").Value = 5940.25
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 5940.25
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 5940.25
$ws.Range("N92").Value = -9684.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5179.9
$ws.Range("I40").Value = 4332.8335
$ws.Range("J40").Value = 6450.5
$ws.Range("K40").Value = 4332.8335
$ws.Range("L40").Value = 6450.5
$ws.Range("M40").Value = -4196.8335
$ws.Range("N40").Value = -6722.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4763171.5
$ws.Range("I122").Value = 9524976
$ws.Range("J122").Value = 1366.6666
$ws.Range("K122").Value = 28574928
$ws.Range("L122").Value = 4099.9998
$ws.Range("M122").Value = -28572478
$ws.Range("N122").Value = -8999.9998
